# Update cryptos list cell values per commit (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.452.80'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.18%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.819.58'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.27%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '314.97'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.46%  '
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.07%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5090'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -4.81%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3949'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -1.44%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08310'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +8.72%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.02%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '41.56'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.78%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.321'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.04%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.05'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.69%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.001'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.15%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.525'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -1.49%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.816.13'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.66%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001148'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +6.66%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '92.53'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +3.22%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06652'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.96%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.78'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.79%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.02%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.123'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.85%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '28.481.92'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.22%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.49'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +3.41%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.271'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +1.76%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '21.29'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +3.09%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '155.98'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.77%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.024.36'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.69%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.72%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '125.69'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +1.54%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.97%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.1093'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -2.13%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.796'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +2.33%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.14%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.07063'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -4.63%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.2225'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.45%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02338'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.05%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.232'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.05%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.867'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.09%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.6301'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.61%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '11.29'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.04%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.16%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.09%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.401'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.57%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.53'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.63%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5921'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +1.42%  '
$ws.Range('B47').Value = 'PancakeSwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.733'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.95%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '125.31'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.41%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.985'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.93%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.185'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -1.38%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06893'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.01%  '
